$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77

# Column A holds a date-like string ("2026/01/26") but must stay plain text,
# matching the existing rows which store dates as literal inline strings
# rather than numeric date values. Force text format before assigning the
# value so Excel does not auto-convert it into a date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2026/01/26"

$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1159

# Match the formatting (centered horizontally/vertically) used by the rest
# of the data rows.
$target = $ws.Range("A" + $row + ":C" + $row)
$target.HorizontalAlignment = -4108
$target.VerticalAlignment = -4108
